$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30228
$ws.Range("H88").Value = 3933.3333
$ws.Range("I88").Value = 4300
$ws.Range("J88").Value = 3750
$ws.Range("K88").Value = 4300
$ws.Range("L88").Value = 3750
$ws.Range("M88").Value = -3894
$ws.Range("N88").Value = -4562
$ws.Range("H91").Value = 3933.3333
$ws.Range("I91").Value = 4300
$ws.Range("J91").Value = 3750
$ws.Range("K91").Value = 4300
$ws.Range("L91").Value = 3750
$ws.Range("M91").Value = -2896
$ws.Range("N91").Value = -6558
$ws.Range("H92").Value = 790.9697
$ws.Range("I92").Value = 329.28
$ws.Range("J92").Value = 2233.75
$ws.Range("K92").Value = 329.28
$ws.Range("L92").Value = 2233.75
$ws.Range("M92").Value = 918.72
$ws.Range("N92").Value = -4729.75
$ws.Range("H93").Value = 40360.6
$ws.Range("J93").Value = 40360.6
$ws.Range("L93").Value = 40360.6
$ws.Range("N93").Value = -45352.6
$ws.Range("H97").Value = 13893.105
$ws.Range("J97").Value = 13893.105
$ws.Range("L97").Value = 41679.315
$ws.Range("N97").Value = -42671.315
$ws.Range("H99").Value = 3033.6667
$ws.Range("I99").Value = 2786
$ws.Range("J99").Value = 3157.5
$ws.Range("K99").Value = 8358
$ws.Range("L99").Value = 9472.5
$ws.Range("M99").Value = -6860
$ws.Range("N99").Value = -12468.5
$ws.Range("H101").Value = 4658.615
$ws.Range("I101").Value = 10441.4
$ws.Range("J101").Value = 1044.375
$ws.Range("K101").Value = 31324.2
$ws.Range("L101").Value = 3133.125
$ws.Range("M101").Value = -29702.2
$ws.Range("N101").Value = -6377.125
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490
$ws.Range("H103").Value = 1017.8
$ws.Range("I103").Value = 900
$ws.Range("J103").Value = 1096.3334
$ws.Range("K103").Value = 2700
$ws.Range("L103").Value = 3289.0002
$ws.Range("M103").Value = -2114
$ws.Range("N103").Value = -4461.0002
$ws.Range("H107").Value = 7026.933
$ws.Range("I107").Value = 7954.231
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 7954.231
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = -6034.231
$ws.Range("N107").Value = -4839.5
$ws.Range("H109").Value = 40093.5
$ws.Range("J109").Value = 40093.5
$ws.Range("L109").Value = 40093.5
$ws.Range("N109").Value = -42867.5
$ws.Range("H111").Value = 2135.6667
$ws.Range("I111").Value = 3250
$ws.Range("J111").Value = 1578.5
$ws.Range("K111").Value = 9750
$ws.Range("L111").Value = 4735.5
$ws.Range("M111").Value = -6683
$ws.Range("N111").Value = -10869.5
$ws.Range("H112").Value = 1201.6111
$ws.Range("I112").Value = 666.6667
$ws.Range("J112").Value = 1308.6
$ws.Range("K112").Value = 2000.0001
$ws.Range("L112").Value = 3925.8
$ws.Range("M112").Value = -892.0001
$ws.Range("N112").Value = -6141.799999999999
$ws.Range("H113").Value = 13698.895
$ws.Range("I113").Value = 41397.8
$ws.Range("J113").Value = 3806.4285
$ws.Range("K113").Value = 41397.8
$ws.Range("L113").Value = 3806.4285
$ws.Range("M113").Value = -38143.8
$ws.Range("N113").Value = -10314.4285
$ws.Range("H137").Value = 478427.47
$ws.Range("I137").Value = 794671.25
$ws.Range("J137").Value = 47185.91
$ws.Range("K137").Value = 2384013.75
$ws.Range("L137").Value = 141557.73
$ws.Range("M137").Value = -2381463.75
$ws.Range("N137").Value = -146657.73
$ws.Range("H141").Value = 2806.56
$ws.Range("I141").Value = 2485.3914
$ws.Range("K141").Value = 7456.174199999999
$ws.Range("M141").Value = -2276.174199999999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9453
$ws.Range("I61").Value = 10878.75
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 10878.75
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -10666.75
$ws.Range("N61").Value = -4174
$ws.Range("H74").Value = 3997.1177
$ws.Range("I74").Value = 748.55
$ws.Range("J74").Value = 8637.929
$ws.Range("K74").Value = 748.55
$ws.Range("L74").Value = 8637.929
$ws.Range("M74").Value = 125.45
$ws.Range("N74").Value = -10385.929
$ws.Range("H77").Value = 3997.1177
$ws.Range("I77").Value = 748.55
$ws.Range("J77").Value = 8637.929
$ws.Range("K77").Value = 3742.75
$ws.Range("L77").Value = 43189.645
$ws.Range("M77").Value = 625.25
$ws.Range("N77").Value = -51925.645
$ws.Range("H132").Value = 3126513.8
$ws.Range("I132").Value = 4167638
$ws.Range("J132").Value = 3140.8
$ws.Range("K132").Value = 12502914
$ws.Range("L132").Value = 9422.400000000001
$ws.Range("M132").Value = -12500384
$ws.Range("N132").Value = -14482.4
$ws.Range("H136").Value = 9453
$ws.Range("I136").Value = 10878.75
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 32636.25
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -30086.25
$ws.Range("N136").Value = -16350

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10117759
$ws.Range("I134").Value = 15898197
$ws.Range("J134").Value = 1994.3334
$ws.Range("K134").Value = 47694591
$ws.Range("L134").Value = 5983.0002
$ws.Range("M134").Value = -47692056
$ws.Range("N134").Value = -11053.0002

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14948.306
$ws.Range("I31").Value = 4411.607
$ws.Range("J31").Value = 51826.75
$ws.Range("K31").Value = 4411.607
$ws.Range("L31").Value = 51826.75
$ws.Range("M31").Value = -4116.607
$ws.Range("N31").Value = -52416.75
$ws.Range("H34").Value = 14948.306
$ws.Range("I34").Value = 4411.607
$ws.Range("J34").Value = 51826.75
$ws.Range("K34").Value = 4411.607
$ws.Range("L34").Value = 51826.75
$ws.Range("M34").Value = -4209.607
$ws.Range("N34").Value = -52230.75
$ws.Range("H58").Value = 15907458
$ws.Range("I58").Value = 71429320
$ws.Range("J58").Value = 44069.715
$ws.Range("K58").Value = 71429320
$ws.Range("L58").Value = 44069.715
$ws.Range("M58").Value = -71429117
$ws.Range("N58").Value = -44475.715
$ws.Range("H132").Value = 22231026
$ws.Range("I132").Value = 41667670
$ws.Range("J132").Value = 17716
$ws.Range("K132").Value = 125003010
$ws.Range("L132").Value = 53148
$ws.Range("M132").Value = -125000480
$ws.Range("N132").Value = -58208
$ws.Range("H134").Value = 10778020
$ws.Range("I134").Value = 15627412
$ws.Range("J134").Value = 4809537
$ws.Range("K134").Value = 46882236
$ws.Range("L134").Value = 14428611
$ws.Range("M134").Value = -46879701
$ws.Range("N134").Value = -14433681
$ws.Range("H136").Value = 15907458
$ws.Range("I136").Value = 71429320
$ws.Range("J136").Value = 44069.715
$ws.Range("K136").Value = 214287960
$ws.Range("L136").Value = 132209.145
$ws.Range("M136").Value = -214285410
$ws.Range("N136").Value = -137309.145

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1386.75
$ws.Range("I98").Value = 1486.6666
$ws.Range("J98").Value = 1326.8
$ws.Range("K98").Value = 4459.9998
$ws.Range("L98").Value = 3980.4
$ws.Range("M98").Value = -2961.9998
$ws.Range("N98").Value = -6976.4
$ws.Range("H141").Value = 4087.5
$ws.Range("I141").Value = 2453.3333
$ws.Range("J141").Value = 8990
$ws.Range("K141").Value = 7359.999899999999
$ws.Range("L141").Value = 26970
$ws.Range("M141").Value = -2179.999899999999
$ws.Range("N141").Value = -37330

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 38466384
$ws.Range("I132").Value = 55556796
$ws.Range("J132").Value = 12953.125
$ws.Range("K132").Value = 166670388
$ws.Range("L132").Value = 38859.375
$ws.Range("M132").Value = -166667858
$ws.Range("N132").Value = -43919.375

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 40912256
$ws.Range("I132").Value = 100001500
$ws.Range("J132").Value = 11367636
$ws.Range("K132").Value = 300004500
$ws.Range("L132").Value = 34102908
$ws.Range("M132").Value = -300001970
$ws.Range("N132").Value = -34107968
$ws.Range("H136").Value = 1157748.4
$ws.Range("I136").Value = 8818.647000000001
$ws.Range("J136").Value = 2552877.5
$ws.Range("K136").Value = 26455.941
$ws.Range("L136").Value = 7658632.5
$ws.Range("M136").Value = -23905.941
$ws.Range("N136").Value = -7663732.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35710820
$ws.Range("I132").Value = 33334602
$ws.Range("J132").Value = 42839480
$ws.Range("K132").Value = 100003806
$ws.Range("L132").Value = 128518440
$ws.Range("M132").Value = -100001276
$ws.Range("N132").Value = -128523500
$ws.Range("H136").Value = 29559064
$ws.Range("I136").Value = 15530341
$ws.Range("J136").Value = 83335830
$ws.Range("K136").Value = 46591023
$ws.Range("L136").Value = 250007490
$ws.Range("M136").Value = -46588473
$ws.Range("N136").Value = -250012590
